$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 1594
$ws.Range("E2").Value = -109
$ws.Range("F2").Value = -109
$ws.Range("G2").Value = -197
$ws.Range("H2").Value = -172
$ws.Range("I2").Value = -161
$ws.Range("J2").Value = -11
$ws.Range("K2").Value = 1865
$ws.Range("L2").Value = 663
$ws.Range("M2").Value = 1202
$ws.Range("N2").Value = 997
$ws.Range("O2").Value = 205
$ws.Range("P2").Value = 170
$ws.Range("Q2").Value = -38
$ws.Range("R2").Value = -122
$ws.Range("S2").Value = 93
$ws.Range("T2").Value = 473
$ws.Range("U2").Value = -511
$ws.Range("V2").Value = 369
$ws.Range("W2").Value = -6.86
$ws.Range("X2").Value = -10.8
$ws.Range("Y2").Value = -14.97
$ws.Range("Z2").Value = -8.94
$ws.Range("AA2").Value = 55.13
$ws.Range("AB2").Value = 439.03
$ws.Range("AC2").Value = -473
$ws.Range("AD2").Value = -3.22
$ws.Range("AE2").Value = 2933
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 34000000

# Row 3 updates
$ws.Range("D3").Value = 1261
$ws.Range("E3").Value = -125
$ws.Range("F3").Value = -125
$ws.Range("G3").Value = -137
$ws.Range("H3").Value = -128
$ws.Range("I3").Value = -116
$ws.Range("J3").Value = -12
$ws.Range("K3").Value = 1702
$ws.Range("L3").Value = 626
$ws.Range("M3").Value = 1077
$ws.Range("N3").Value = 884
$ws.Range("O3").Value = 192
$ws.Range("P3").Value = 170
$ws.Range("Q3").Value = 50
$ws.Range("R3").Value = -50
$ws.Range("S3").Value = -5
$ws.Range("T3").Value = 84
$ws.Range("U3").Value = -35
$ws.Range("V3").Value = 364
$ws.Range("W3").Value = -9.88
$ws.Range("X3").Value = -10.15
$ws.Range("Y3").Value = -12.28
$ws.Range("Z3").Value = -7.17
$ws.Range("AA3").Value = 58.14
$ws.Range("AB3").Value = 370.79
$ws.Range("AC3").Value = -340
$ws.Range("AD3").Value = -5.25
$ws.Range("AE3").Value = 2601
$ws.Range("AF3").Value = 0.69
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 34000000

# Row 4 updates
$ws.Range("D4").Value = 1136
$ws.Range("E4").Value = -64
$ws.Range("F4").Value = -64
$ws.Range("G4").Value = -109
$ws.Range("H4").Value = -107
$ws.Range("I4").Value = -70
$ws.Range("J4").Value = -37
$ws.Range("K4").Value = 1566
$ws.Range("L4").Value = 603
$ws.Range("M4").Value = 963
$ws.Range("N4").Value = 845
$ws.Range("O4").Value = 117
$ws.Range("P4").Value = 170
$ws.Range("Q4").Value = 85
$ws.Range("R4").Value = -20
$ws.Range("S4").Value = -51
$ws.Range("T4").Value = 37
$ws.Range("U4").Value = 48
$ws.Range("V4").Value = 313
$ws.Range("W4").Value = -5.59
$ws.Range("X4").Value = -9.46
$ws.Range("Y4").Value = -8.09
$ws.Range("Z4").Value = -6.57
$ws.Range("AA4").Value = 62.67
$ws.Range("AB4").Value = 337.44
$ws.Range("AC4").Value = -206
$ws.Range("AD4").Value = -7.63
$ws.Range("AE4").Value = 2488
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 34000000

# Row 5 updates
$ws.Range("D5").Value = 1004
$ws.Range("E5").Value = -138
$ws.Range("F5").Value = -138
$ws.Range("G5").Value = -298
$ws.Range("H5").Value = -318
$ws.Range("I5").Value = -215
$ws.Range("J5").Value = -103
$ws.Range("K5").Value = 1127
$ws.Range("L5").Value = 582
$ws.Range("M5").Value = 545
$ws.Range("N5").Value = 531
$ws.Range("O5").Value = 14
$ws.Range("P5").Value = 170
$ws.Range("Q5").Value = -57
$ws.Range("R5").Value = 77
$ws.Range("S5").Value = 46
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = -85
$ws.Range("V5").Value = 350
$ws.Range("W5").Value = -13.74
$ws.Range("X5").Value = -31.69
$ws.Range("Y5").Value = -31.2
$ws.Range("Z5").Value = -23.61
$ws.Range("AA5").Value = 106.8
$ws.Range("AB5").Value = 212.46
$ws.Range("AC5").Value = -631
$ws.Range("AD5").Value = -1.81
$ws.Range("AE5").Value = 1562
$ws.Range("AF5").Value = 0.73
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 34000000
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6 updates
$ws.Range("D6").Value = 1029
$ws.Range("E6").Value = -63
$ws.Range("F6").Value = -63
$ws.Range("G6").Value = -77
$ws.Range("H6").Value = -78
$ws.Range("I6").Value = -69
$ws.Range("K6").Value = 1052
$ws.Range("L6").Value = 589
$ws.Range("M6").Value = 463
$ws.Range("N6").Value = 459
$ws.Range("P6").Value = 170
$ws.Range("Q6").Value = -38
$ws.Range("R6").Value = -36
$ws.Range("S6").Value = 14
$ws.Range("T6").Value = 44
$ws.Range("U6").Value = -82
$ws.Range("V6").Value = 352
$ws.Range("W6").Value = -6.09
$ws.Range("X6").Value = -7.59
$ws.Range("Y6").Value = -13.92
$ws.Range("Z6").Value = -7.17
$ws.Range("AA6").Value = 127.04
$ws.Range("AB6").Value = 169.54
$ws.Range("AC6").Value = -203
$ws.Range("AD6").Value = -4.09
$ws.Range("AE6").Value = 1349
$ws.Range("AF6").Value = 0.61
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 34000000
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7: clear all data columns, keep A/B/C only
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: clear all data columns, keep A/B/C only
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: clear all data columns, keep A/B/C only
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

